# Automatische test-sync: 2025-08-06 20:39:50
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append new log row (row 19) -----------------------------------------
$newRow = 19
$ws.Cells.Item($newRow, 1).Value = "Weten jullie al iets over mijn retour?"
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Testmail #1: Weten jullie al iets over mijn retour?"
$ws.Cells.Item($newRow, 4).Value = "Retour / Terugbetaling"
$ws.Cells.Item($newRow, 5).Value = "Beste klant,`nBedankt voor je e-mail. Om je vraag over je retour te beantwoorden, heb ik wat meer informatie nodig. Zou je ons je ordernummer of het track & trace-nummer van de retourzending kunnen doorgeven? Op die manier kunnen we de status controleren en je verder helpen.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Cells.Item($newRow, 6).Value = "2025-08-06 20:39:04"
$ws.Cells.Item($newRow, 7).Value = "Ja"
$ws.Cells.Item($newRow, 8).Value = "Nee"
$ws.Cells.Item($newRow, 9).Value = "Ja"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Re-fit the new row's height (clears the explicit/custom height that gets
# pinned by the multi-line write, keeping the row on the default height
# like every other row in the sheet).
$ws.Rows.Item($newRow).AutoFit()

# --- Extend conditional formatting ranges to include the new row ---------
foreach ($col in @("D","G","H","I","J")) {
    $fcs = $ws.Range($col + "2").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($col + "2:" + $col + $newRow))
    }
}

# --- Update Dashboard summary (re-sorted by count, descending) -----------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(4, 2).Value = 3
$dash.Cells.Item(5, 1).Value = "Klantenservice / Contact"
$dash.Cells.Item(5, 2).Value = 2
$dash.Cells.Item(6, 1).Value = "Klantenservice / Opvolging"
$dash.Cells.Item(6, 2).Value = 2
